$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# hardmode_percent (column E) was stored as a fraction (e.g. 0.0169...)
# but should actually be expressed as a percentage (e.g. 1.69...).
# Multiply every data row's value in column E by 100.

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 354 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value2 = $val * 100
    }
}
